$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 54, pushing the existing rows 54-55 down to 55-56.
$ws.Rows.Item(54).Insert()

# Row 53 already carries the exact border/number-format style pattern
# (A=1, B=16, C=18, D=17, E=18, F=18, G=16) that the new "Task - Only
# participants can add tags" row needs, so copy its formatting down into
# the freshly inserted row 54.
$ws.Range("A53:G53").Copy()
$ws.Range("A54:G54").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New row is taller (two lines of label text) than its donor row.
$ws.Rows.Item(54).RowHeight = 45

# Populate the new rule row.
$ws.Range("B54").Value = "Task – Only participants can add tags"
$ws.Range("C54").Value = "TASK"
$ws.Range("G54").Value = "grant addTag to assignee, co-owner, supervisor, owning group, approver, collaborator, reader"

# Clear any stray formula-esque bits the paste might have introduced; D54
# must stay blank (style 17 only, same as D53 before population).
$ws.Range("D54").ClearContents()

# Move the visible selection to mirror the author's final cursor position.
$ws.Range("C54").Select()
